$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.544.78'
$ws.Range('E2').Value = '  -0.66%  '
$ws.Range('D3').Value = '1.911.88'
$ws.Range('E3').Value = '  -1.36%  '
$ws.Range('E4').Value = '  -0.17%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '239.42'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  -1.09%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '1.001'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  -0.18%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.4783'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  -1.93%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '0.2842'
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  -2.70%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.06690'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  -2.54%  '
$ws.Range('E10').Value = '  -3.95%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '101.21'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  -3.74%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.913.10'
$ws.Range('E12').Value = '  -1.33%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '0.07685'
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  -0.86%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '5.209'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  -1.85%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '0.6690'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  -3.81%  '
$ws.Range('D16').Value = '30.528.96'
$ws.Range('E16').Value = '  -0.78%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '256.44'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  -6.73%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '1.000'
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  -0.16%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '0.000007461'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  -3.05%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '12.64'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  -3.65%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '5.381'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  -1.05%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '1.001'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  -0.25%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '6.290'
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  -2.49%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '9.324'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  -3.73%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '166.77'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  -0.44%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '19.16'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  -1.93%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '2.057'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  -4.73%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '4.743'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  +4.46%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '1.382'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  -0.57%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '0.1007'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  -3.03%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '1.512'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  -2.43%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '4.236'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  -2.51%  '
$ws.Range('E33').Value = '  -2.53%  '
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '0.7274'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  -2.54%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '1.108'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  -3.85%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '0.9997'
$cell.Style = 'Normal'
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '2.709'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  -0.74%  '
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '0.01910'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  -3.78%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '2.616'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  -1.87%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '6.242'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  -3.00%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '74.53'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  -3.68%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '1.965'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  -5.45%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '0.8610'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  -3.73%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '105.40'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  -2.43%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '1.000'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  +0.15%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '0.4229'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  -3.81%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '7.345'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  -5.07%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '0.1199'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  -3.00%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '34.69'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  -2.80%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '907.19'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  -9.43%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '8.782'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  -4.37%  '
